$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Document Type")

# --- Fixed French Credit Note categories (rows 278-279): Invoice -> Credit Note ---
$ws.Cells.Item(278, 13).Value = "Credit Note"
$ws.Cells.Item(279, 13).Value = "Credit Note"

# --- Added AE PINT DocTypes: four new rows (280-283) ---

# Column A (Profile name) - written first so new shared strings land in this order
$ws.Cells.Item(280, 1).Value = "AE PINT Invoice v1.0"
$ws.Cells.Item(281, 1).Value = "AE PINT Credit Note v1.0"
$ws.Cells.Item(282, 1).Value = "AE PINT Self-Billing Invoice v1.0"
$ws.Cells.Item(283, 1).Value = "AE PINT Self-Billing Credit Note v1.0"

# Column C (Peppol Document Type Identifier Value)
$ws.Cells.Item(280, 3).Value = "urn:oasis:names:specification:ubl:schema:xsd:Invoice-2::Invoice##urn:peppol:pint:billing-1@ae-1::2.1"
$ws.Cells.Item(281, 3).Value = "urn:oasis:names:specification:ubl:schema:xsd:CreditNote-2::CreditNote##urn:peppol:pint:billing-1@ae-1::2.1"
$ws.Cells.Item(282, 3).Value = "urn:oasis:names:specification:ubl:schema:xsd:Invoice-2::Invoice##urn:peppol:pint:selfbilling-1@ae-1::2.1"
$ws.Cells.Item(283, 3).Value = "urn:oasis:names:specification:ubl:schema:xsd:CreditNote-2::CreditNote##urn:peppol:pint:selfbilling-1@ae-1::2.1"

# Column H (Comment) - TICC-373, shared across all four new rows
$ws.Cells.Item(280, 8).Value = "TICC-373"
$ws.Cells.Item(281, 8).Value = "TICC-373"
$ws.Cells.Item(282, 8).Value = "TICC-373"
$ws.Cells.Item(283, 8).Value = "TICC-373"

# Column B (Peppol Document Type Identifier Scheme)
$ws.Cells.Item(280, 2).Value = "peppol-doctype-wildcard"
$ws.Cells.Item(281, 2).Value = "peppol-doctype-wildcard"
$ws.Cells.Item(282, 2).Value = "peppol-doctype-wildcard"
$ws.Cells.Item(283, 2).Value = "peppol-doctype-wildcard"

# Column D (Initial release)
$ws.Cells.Item(280, 4).Value = "9.1"
$ws.Cells.Item(281, 4).Value = "9.1"
$ws.Cells.Item(282, 4).Value = "9.1"
$ws.Cells.Item(283, 4).Value = "9.1"

# Column E (State)
$ws.Cells.Item(280, 5).Value = "active"
$ws.Cells.Item(281, 5).Value = "active"
$ws.Cells.Item(282, 5).Value = "active"
$ws.Cells.Item(283, 5).Value = "active"

# Column I (Abstract?)
$ws.Cells.Item(280, 9).Value = $False
$ws.Cells.Item(281, 9).Value = $False
$ws.Cells.Item(282, 9).Value = $False
$ws.Cells.Item(283, 9).Value = $False

# Column J (Issued by OpenPeppol?)
$ws.Cells.Item(280, 10).Formula = "=TRUE"
$ws.Cells.Item(281, 10).Formula = "=TRUE"
$ws.Cells.Item(282, 10).Formula = "=TRUE"
$ws.Cells.Item(283, 10).Formula = "=TRUE"

# Column K (BIS version)
$ws.Cells.Item(280, 11).Value = 3
$ws.Cells.Item(281, 11).Value = 3
$ws.Cells.Item(282, 11).Value = 3
$ws.Cells.Item(283, 11).Value = 3

# Column L (Domain Community)
$ws.Cells.Item(280, 12).Value = "POAC"
$ws.Cells.Item(281, 12).Value = "POAC"
$ws.Cells.Item(282, 12).Value = "POAC"
$ws.Cells.Item(283, 12).Value = "POAC"

# Column M (Category)
$ws.Cells.Item(280, 13).Value = "Invoice"
$ws.Cells.Item(281, 13).Value = "Credit Note"
$ws.Cells.Item(282, 13).Value = "Invoice"
$ws.Cells.Item(283, 13).Value = "Credit Note"

# Column N (Associated Process/Profile Identifier(s))
$ws.Cells.Item(280, 14).Value = "cenbii-procid-ubl::urn:peppol:bis:billing"
$ws.Cells.Item(281, 14).Value = "cenbii-procid-ubl::urn:peppol:bis:billing"
$ws.Cells.Item(282, 14).Value = "cenbii-procid-ubl::urn:peppol:bis:selfbilling"
$ws.Cells.Item(283, 14).Value = "cenbii-procid-ubl::urn:peppol:bis:selfbilling"

# Apply cell formatting to match the equivalent existing "MY PINT" block (rows
# 246-249), which uses the exact same style pattern for this kind of doc-type row.
# Done last (after values) and column by column so the paste doesn't disturb the
# values/formulas just written and doesn't create stray cells in columns F/G.
$ws.Range("A246:A249").Copy()
$ws.Range("A280:A283").PasteSpecial(-4122)
$excel.CutCopyMode = $False

$ws.Range("B246:B249").Copy()
$ws.Range("B280:B283").PasteSpecial(-4122)
$excel.CutCopyMode = $False

$ws.Range("C246:C249").Copy()
$ws.Range("C280:C283").PasteSpecial(-4122)
$excel.CutCopyMode = $False

$ws.Range("D246:D249").Copy()
$ws.Range("D280:D283").PasteSpecial(-4122)
$excel.CutCopyMode = $False

$ws.Range("E246:E249").Copy()
$ws.Range("E280:E283").PasteSpecial(-4122)
$excel.CutCopyMode = $False

$ws.Range("H246:N249").Copy()
$ws.Range("H280:N283").PasteSpecial(-4122)
$excel.CutCopyMode = $False

# Rows 281 & 283 wrap onto two lines (like their 246-249 template counterparts), so
# match the taller row height; 280 & 282 stay at the default row height.
$ws.Rows.Item(281).RowHeight = 30
$ws.Rows.Item(283).RowHeight = 30

$ws.Range("A283").Select()
